# Adding attempt to load data from MS CSV
# Appends a new day of data (2020-04-04) as row 41 to the "Confirmados"
# and "Mortes" sheets, mirroring the existing layout (Data in col A,
# one column per UF from B..AB).

$wb = $excel.ActiveWorkbook

$wsConfirmados = $wb.Worksheets.Item("Confirmados")
$wsMortes      = $wb.Worksheets.Item("Mortes")

$newDate = "2020-04-04"

# Values for each UF column (B..AB), in the same order as the header row.
$confirmados = @(46,23,28,311,332,730,454,153,103,88,56,62,430,80,32,395,176,22,1246,212,410,11,37,334,4466,27,14)
$mortes      = @(0,2,1,12,7,22,7,5,2,1,1,1,6,1,2,6,14,4,58,5,6,1,1,5,260,2,0)

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB")

$rowIndex = 41

# -- Confirmados sheet --
# Force the date to be written as plain text (matching the existing
# "Data" column, which stores dates as text strings, not date serials),
# then restore the cell style to the plain/default one used elsewhere
# in the column so no stray number-format style sticks around.
$wsConfirmados.Range("A$rowIndex").NumberFormat = "@"
$wsConfirmados.Range("A$rowIndex").Value = $newDate
$wsConfirmados.Range("A$rowIndex").Style = $wsConfirmados.Range("A2").Style
for ($i = 0; $i -lt $cols.Length; $i++) {
    $wsConfirmados.Range("$($cols[$i])$rowIndex").Value = $confirmados[$i]
}

# -- Mortes sheet --
$wsMortes.Range("A$rowIndex").NumberFormat = "@"
$wsMortes.Range("A$rowIndex").Value = $newDate
$wsMortes.Range("A$rowIndex").Style = $wsMortes.Range("A2").Style
for ($i = 0; $i -lt $cols.Length; $i++) {
    $wsMortes.Range("$($cols[$i])$rowIndex").Value = $mortes[$i]
}
